$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "24.005.53"
$ws.Range("E2").Value = "  -1.90%  "
Set-TextValue "D3" "1.651.82"
$ws.Range("E3").Value = "  -1.02%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue "D5" "310.16"
$ws.Range("E6").Value = "  +0.11%  "
Set-TextValue "D7" "0.3910"
$ws.Range("E7").Value = "  -1.49%  "
Set-TextValue "D8" "0.3815"
$ws.Range("E8").Value = "  -2.70%  "
Set-TextValue "D9" "52.29"
$ws.Range("E9").Value = "  +0.63%  "
Set-TextValue "D10" "1.350"
$ws.Range("E10").Value = "  -4.26%  "
Set-TextValue "D11" "1.002"
$ws.Range("E11").Value = "  +0.25%  "
Set-TextValue "D12" "0.08457"
$ws.Range("E12").Value = "  -1.77%  "
Set-TextValue "D13" "23.92"
$ws.Range("E13").Value = "  -2.36%  "
Set-TextValue "D14" "7.074"
$ws.Range("E14").Value = "  -3.87%  "
Set-TextValue "D15" "8.016"
$ws.Range("E15").Value = "  +1.20%  "
Set-TextValue "D16" "0.00001311"
$ws.Range("E16").Value = "  -2.91%  "
Set-TextValue "D17" "1.650.82"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  -1.39%  "
Set-TextValue "D19" "0.07010"
$ws.Range("E19").Value = "  +0.35%  "
Set-TextValue "D20" "19.73"
Set-TextValue "D21" "6.986"
$ws.Range("E21").Value = "  -0.55%  "
Set-TextValue "D23" "13.81"
$ws.Range("E23").Value = "  +0.15%  "
Set-TextValue "D24" "23.991.20"
$ws.Range("E24").Value = "  -1.94%  "
Set-TextValue "D25" "2.443"
$ws.Range("E25").Value = "  +0.60%  "
Set-TextValue "D26" "2.953"
$ws.Range("E26").Value = "  -2.67%  "
Set-TextValue "D27" "22.10"
$ws.Range("E27").Value = "  -1.99%  "
Set-TextValue "D28" "152.95"
$ws.Range("E28").Value = "  -2.90%  "
Set-TextValue "D29" "5.408"
$ws.Range("E29").Value = "  -0.53%  "
Set-TextValue "D30" "138.21"
$ws.Range("E30").Value = "  -3.41%  "
Set-TextValue "D31" "7.945"
$ws.Range("E31").Value = "  -2.53%  "
Set-TextValue "D32" "2.507"
$ws.Range("E32").Value = "  -1.41%  "
Set-TextValue "D33" "1.830.58"
$ws.Range("E33").Value = "  -0.98%  "
Set-TextValue "D34" "1.021"
$ws.Range("E34").Value = "  -4.23%  "
Set-TextValue "D35" "0.08060"
$ws.Range("E35").Value = "  -2.72%  "
Set-TextValue "D36" "6.754"
$ws.Range("E36").Value = "  -1.25%  "
Set-TextValue "D37" "0.02926"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  -3.67%  "
Set-TextValue "D39" "0.2680"
$ws.Range("E39").Value = "  -3.18%  "
Set-TextValue "D40" "0.09101"
$ws.Range("E40").Value = "  -1.83%  "
Set-TextValue "D41" "0.7619"
$ws.Range("E41").Value = "  -1.96%  "
Set-TextValue "D42" "13.43"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  -2.67%  "
Set-TextValue "D45" "0.6980"
$ws.Range("E45").Value = "  -2.39%  "
Set-TextValue "D46" "2.458"
$ws.Range("E46").Value = "  -3.40%  "
Set-TextValue "D47" "4.102"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +0.08%  "
Set-TextValue "D49" "0.08323"
$ws.Range("E49").Value = "  -1.57%  "
Set-TextValue "D50" "134.97"
$ws.Range("E50").Value = "  -1.36%  "
Set-TextValue "D51" "1.234"
$ws.Range("E51").Value = "  -3.81%  "
